$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each updated cell is forced to Text format first so that numeric-looking
# strings (prices like "312.99", "1.00", percentages, etc.) are preserved
# verbatim instead of being auto-converted to floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.359.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.468.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.57%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.05"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.850.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.472.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.67"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.313.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0923"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.15"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.74"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.06%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.16"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.55"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.59"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.62"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.82%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0757"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.19%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.97"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.19%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.34"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.990.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "70.01"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.20"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.180"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.61%  "

Write-Host "Applied 94 cell updates"
